$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDM")
$ws.Rows("21:40").Delete()
Write-Host $ws.Range("B20").Value2
Write-Host $ws.Range("B21").Value2
Write-Host $ws.Cells.Item(41,2).Value2
